# aggiunti autonomia, prezzo, airbag
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row for "carrozzeria" right after "modello" (before old "anno_fabbricazione") ---
$ws.Rows("11:11").Insert()
$ws.Range("B11").Value = "carrozzeria"
$ws.Range("C11").Value = "varchar(40)"
$ws.Range("D11").Value = "notnull"
$ws.Range("E11").Value = "es. berlina, SUV ecc."

# --- 2) Insert a new row for "prezzo" right after "chilometri" (before old "alimentazione") ---
$ws.Rows("14:14").Insert()
$ws.Range("B14").Value = "prezzo"
$ws.Range("C14").Value = "int"
$ws.Range("D14").Value = "notnull"

# --- 3) Insert a new row for "autonomia" right after "consumo" (before old "emissioni") ---
$ws.Rows("17:17").Insert()
$ws.Range("B17").Value = "autonomia"
$ws.Range("C17").Value = "smallint"
$ws.Range("E17").Value = "km (se alimentazione elettrica)"

# Update "consumo" row (now row 16): type grows to decimal(4,1), "notnull" attribute removed,
# description clarified as litri/100km (se alimentazione a combustibile)
$ws.Range("C16").Value = "decimal(4, 1)"
$ws.Range("D16").ClearContents()
$ws.Range("E16").Value = "litri/100km (se alimentazione a combustibile)"

# --- 4) Append a new row for "airbag" at the end of the table ---
$ws.Range("B28").Value = "airbag"
$ws.Range("C28").Value = "tinyint"
$ws.Range("D28").Value = "notnull"

# Column E got wider to fit the longer descriptions
$ws.Columns("E:E").ColumnWidth = 41.45

# Selection moved to D16 in the saved file
$ws.Range("D16").Select()
